# Adds the first set of Thalmor cards to the "Units & Decks" workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A194").Value = '##########################################################################################################'
$ws.Range("A200").Value = 'Warriors'
$ws.Range("A201").Value = 'No.'
$ws.Range("B201").Value = 'Name'
$ws.Range("C201").Value = 'Strength'
$ws.Range("D201").Value = 'Subtype'
$ws.Range("E201").Value = 'Created'
$ws.Range("F201").Value = 'Hero'
$ws.Range("G201").Value = 'Race'
$ws.Range("H201").Value = 'Ability'
$ws.Range("A215").Value = 'Mages'
$ws.Range("A216").Value = 'No.'
$ws.Range("B216").Value = 'Name'
$ws.Range("C216").Value = 'Strength'
$ws.Range("D216").Value = 'Subtype'
$ws.Range("E216").Value = 'Created'
$ws.Range("F216").Value = 'Hero'
$ws.Range("G216").Value = 'Race'
$ws.Range("H216").Value = 'Ability'
$ws.Range("A231").Value = 'Shadow'
$ws.Range("A232").Value = 'No.'
$ws.Range("B232").Value = 'Name'
$ws.Range("C232").Value = 'Strength'
$ws.Range("D232").Value = 'Subtype'
$ws.Range("E232").Value = 'Created'
$ws.Range("F232").Value = 'Hero'
$ws.Range("G232").Value = 'Race'
$ws.Range("H232").Value = 'Ability'
$ws.Range("B196").Value = 'Deck: Thalmor'
$ws.Range("B217").Value = 'Thalmor Wizard'
$ws.Range("A217").Value = 1
$ws.Range("C217").Value = 2
$ws.Range("B233").Value = 'Thalmor Archer'
$ws.Range("A233").Value = 1
$ws.Range("C233").Value = 2
$ws.Range("E233").Value = 'yes'
$ws.Range("B202").Value = 'Thalmor Soldier'
$ws.Range("A202").Value = 1
$ws.Range("C202").Value = 1
$ws.Range("E202").Value = 'yes'
$ws.Range("B234").Value = 'Northwatch Archer'
$ws.Range("A234").Value = 2
$ws.Range("C234").Value = 3
$ws.Range("E234").Value = 'yes'
$ws.Range("B218").Value = 'Northwatch Mage'
$ws.Range("A218").Value = 2
$ws.Range("C218").Value = 3
$ws.Range("B203").Value = 'Northwatch Guard'
$ws.Range("A203").Value = 2
$ws.Range("C203").Value = 2
$ws.Range("E203").Value = 'yes'
$ws.Range("B219").Value = 'Northwatch Interogater'
$ws.Range("A219").Value = 3
$ws.Range("C219").Value = 4
$ws.Range("B235").Value = 'Shavari'
$ws.Range("A235").Value = 3
$ws.Range("C235").Value = 6
$ws.Range("D235").Value = 'Spy'
$ws.Range("E235").Value = 'yes'
$ws.Range("D205").Value = 'Spy'
$ws.Range("A205").Value = 4
$ws.Range("C205").Value = 4
$ws.Range("E205").Value = 'yes'
$ws.Range("G235").Value = 'Khajiit'
$ws.Range("B220").Value = 'Ancano'
$ws.Range("A220").Value = 4
$ws.Range("C220").Value = 0
$ws.Range("D220").Value = 'Spy'
$ws.Range("E220").Value = 'yes'
$ws.Range("F220").Value = 'yes'
$ws.Range("G220").Value = 'altmer'
$ws.Range("B204").Value = 'Thalmor Agent'
$ws.Range("A204").Value = 3
$ws.Range("C204").Value = 3
$ws.Range("E204").Value = 'yes'
$ws.Range("B198").Value = 'Elenwen a leader card'
$ws.Range("B221").Value = 'Estormo'
$ws.Range("A221").Value = 5
$ws.Range("C221").Value = 6
$ws.Range("D221").Value = 'Spy'
$ws.Range("E221").Value = 'yes'
$ws.Range("B236").Value = 'Lorcalin'
$ws.Range("A236").Value = 4
$ws.Range("C236").Value = 8
$ws.Range("D236").Value = 'Spy'
$ws.Range("B205").Value = 'Ondolemar'
$ws.Range("B206").Value = 'Justicar'
$ws.Range("A206").Value = 5
$ws.Range("C206").Value = 5
$ws.Range("B207").Value = 'bound blade assassin'
$ws.Range("A207").Value = 6
$ws.Range("C207").Value = 3
$ws.Range("B208").Value = 'warrior'
$ws.Range("A208").Value = 7
$ws.Range("C208").Value = 4
$ws.Range("B223").Value = 'spellsword 1'
$ws.Range("A223").Value = 7
$ws.Range("C223").Value = 4
$ws.Range("E223").Value = 'yes'
$ws.Range("B224").Value = 'spellsword 2'
$ws.Range("A224").Value = 8
$ws.Range("C224").Value = 6
$ws.Range("E224").Value = 'yes'
$ws.Range("B209").Value = 'rulindil'
$ws.Range("A209").Value = 8
$ws.Range("C209").Value = 7
$ws.Range("E209").Value = 'yes'
$ws.Range("A210").Value = 9
$ws.Range("A211").Value = 10
$ws.Range("A212").Value = 11
$ws.Range("A213").Value = 12
$ws.Range("A222").Value = 6
$ws.Range("A237").Value = 5
$ws.Range("A238").Value = 6
$ws.Range("A239").Value = 7
$ws.Range("A240").Value = 8

# Update the view to match where the author left off editing
$ws.Application.ActiveWindow.ScrollRow = 201
$ws.Range("F212").Select()
